{"js": "// Add a new \"Brazilian Portuguese (Adult Baseline)\" list item right after\n// the existing \"Spanish\" list item in the \"v0.3 translations:\" list.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the \"Spanish\" list entry under the \"v0.3 translations:\" heading.\nlet spanishParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Spanish\") {\n    spanishParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!spanishParagraph) {\n  throw new Error('Could not find the \"Spanish\" paragraph to anchor the new entry.');\n}\n\n// Insert a sibling paragraph after it, inheriting the same list\n// formatting (ListParagraph style / numId 10), and set its text.\nspanishParagraph.insertParagraph(\"Brazilian Portuguese (Adult Baseline)\", \"After\");\n\nawait context.sync();\n", "ps1": "# Add a new \"Brazilian Portuguese (Adult Baseline)\" list item right after\n# the existing \"Spanish\" list item in the \"v0.3 translations:\" list.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Spanish\")\n\nif ($found) {\n    $spanishPara = $rng.Paragraphs.Item(1)\n    $spanishIndex = $spanishPara.Index\n\n    # Insert a new sibling paragraph after \"Spanish\"; it inherits the same\n    # list formatting (ListParagraph style / numId 10).\n    $spanishPara.Range.InsertParagraphAfter()\n\n    # Re-fetch the freshly created paragraph by its position in the\n    # document's paragraph collection and set its text.\n    $newPara = $d.Paragraphs.Item($spanishIndex + 1)\n    $newPara.Range.Text = \"Brazilian Portuguese (Adult Baseline)\"\n}\n"}
